$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''26.609.22'
$ws.Range("E2").Value = '  +0.66%  '
$ws.Range("D3").Value = '''1.744.53'
$ws.Range("E3").Value = '  +1.15%  '
$ws.Range("D4").Value = '''0.9998'
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '''246.67'
$ws.Range("E5").Value = '  +0.89%  '
$ws.Range("D6").Value = '''0.9999'
$ws.Range("E6").Value = '  -0.08%  '
$ws.Range("D7").Value = '''0.4920'
$ws.Range("E7").Value = '  +2.64%  '
$ws.Range("D8").Value = '''0.2680'
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("D9").Value = '''0.06321'
$ws.Range("E9").Value = '  +1.65%  '
$ws.Range("D10").Value = '''1.744.25'
$ws.Range("E10").Value = '  +1.03%  '
$ws.Range("D11").Value = '''0.07048'
$ws.Range("E11").Value = '  -0.89%  '
$ws.Range("E12").Value = '  +0.36%  '
$ws.Range("D13").Value = '''0.6158'
$ws.Range("D14").Value = '''4.584'
$ws.Range("D15").Value = '''78.17'
$ws.Range("E15").Value = '  +1.34%  '
$ws.Range("E16").Value = '  -0.06%  '
$ws.Range("D17").Value = '''26.621.64'
$ws.Range("E17").Value = '  +0.62%  '
$ws.Range("D18").Value = '''0.000007327'
$ws.Range("E18").Value = '  +5.64%  '
$ws.Range("D19").Value = '''1.0000'
$ws.Range("E19").Value = '  -0.10%  '
$ws.Range("D20").Value = '''11.59'
$ws.Range("E20").Value = '  -0.94%  '
$ws.Range("D21").Value = '''1.967.59'
$ws.Range("E21").Value = '  +0.63%  '
$ws.Range("D22").Value = '''4.575'
$ws.Range("D23").Value = '''8.729'
$ws.Range("E23").Value = '  -1.88%  '
$ws.Range("D24").Value = '''5.267'
$ws.Range("E24").Value = '  -0.80%  '
$ws.Range("D25").Value = '''139.74'
$ws.Range("E25").Value = '  +2.59%  '
$ws.Range("D26").Value = '''15.46'
$ws.Range("E26").Value = '  +0.80%  '
$ws.Range("D27").Value = '''1.427'
$ws.Range("E27").Value = '  +1.08%  '
$ws.Range("D28").Value = '''1.765'
$ws.Range("E28").Value = '  -1.62%  '
$ws.Range("D29").Value = '''107.66'
$ws.Range("E29").Value = '  +0.95%  '
$ws.Range("D30").Value = '''4.048'
$ws.Range("E30").Value = '  +2.01%  '
$ws.Range("D31").Value = '''0.08042'
$ws.Range("E31").Value = '  +0.11%  '
$ws.Range("D32").Value = '''3.742'
$ws.Range("E32").Value = '  +0.20%  '
$ws.Range("D33").Value = '''0.04626'
$ws.Range("E33").Value = '  +1.81%  '
$ws.Range("D34").Value = '''0.9993'
$ws.Range("E34").Value = '  -0.07%  '
$ws.Range("E35").Value = '  -0.30%  '
$ws.Range("D36").Value = '''1.019'
$ws.Range("E36").Value = '  +3.23%  '
$ws.Range("D37").Value = '''0.6398'
$ws.Range("E37").Value = '  +0.55%  '
$ws.Range("D38").Value = '''2.080'
$ws.Range("E38").Value = '  +4.46%  '
$ws.Range("D39").Value = '''0.8998'
$ws.Range("E39").Value = '  -3.82%  '
$ws.Range("D40").Value = '''2.427'
$ws.Range("E40").Value = '  +1.03%  '
$ws.Range("D41").Value = '''1.003'
$ws.Range("E41").Value = '  -0.23%  '
$ws.Range("D42").Value = '''0.01506'
$ws.Range("E42").Value = '  +0.58%  '
$ws.Range("D43").Value = '''102.00'
$ws.Range("E43").Value = '  -5.04%  '
$ws.Range("D44").Value = '''5.430'
$ws.Range("E44").Value = '  -3.58%  '
$ws.Range("D45").Value = '''0.3926'
$ws.Range("E45").Value = '  +0.49%  '
$ws.Range("D46").Value = '''6.882'
$ws.Range("E46").Value = '  -1.17%  '
$ws.Range("D47").Value = '''0.1184'
$ws.Range("E47").Value = '  -0.62%  '
$ws.Range("D48").Value = '''0.05397'
$ws.Range("E48").Value = '  +1.54%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '''7.841'
$ws.Range("E49").Value = '  -0.03%  '
$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").Value = '''30.58'
$ws.Range("E50").Value = '  -1.05%  '
$ws.Range("D51").Value = '''1.263'
$ws.Range("E51").Value = '  -0.35%  '
